$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the staff name to include the "Mr." honorific prefix (multi-line, trailing newline retained)
$ws.Range("A2").Value = "Mr. MUTHIAH M M`n"

# Reset explicit per-cell style on the other row-2 cells back to the default "Normal" style
$ws.Range("B2:D2").Style = "Normal"
$ws.Range("I2:J2").Style = "Normal"

# Remove the unused blank cells (Research Gate / Orchid / Publon / Scopus were empty)
$ws.Range("E2:H2").Clear()

# Select A2 as the active cell (matches final selection in the sheet view)
$ws.Range("A2").Select()
